$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.3997546666666666
$ws.Range("H2").Value = 1.199264
$ws.Range("I2").Value = 0.1320462084214824
$ws.Range("J2").Value = 0.1320462084214824
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 6.298390483068444
$ws.Range("R2").Value = 56.685514347616
$ws.Range("S2").Value = 0.04252623294978246
$ws.Range("T2").Value = 0.04252623294978245
$ws.Range("G3").Value = 0.3997546666666666
$ws.Range("H3").Value = 1.199264
$ws.Range("I3").Value = 0.1320462084214824
$ws.Range("J3").Value = 0.1320462084214824
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("Q3").Value = 10.77492625612089
$ws.Range("R3").Value = 96.97433630508799
$ws.Range("S3").Value = 0.07275144740808298
$ws.Range("T3").Value = 0.07275144740808295
$ws.Range("G4").Value = 0.3997546666666666
$ws.Range("H4").Value = 1.199264
$ws.Range("I4").Value = 0.1320462084214824
$ws.Range("J4").Value = 0.1320462084214824
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 2.483519706428444
$ws.Range("R4").Value = 22.351677357856
$ws.Range("S4").Value = 0.01676852806361699
$ws.Range("T4").Value = 0.01676852806361698
$ws.Range("I5").Value = 0.6840925621829359
$ws.Range("J5").Value = 0.684092562182936
$ws.Range("M5").Value = 15.75563966666667
$ws.Range("N5").Value = 47.266919
$ws.Range("O5").Value = 0.3220556913988901
$ws.Range("P5").Value = 0.32205569139889
$ws.Range("Q5").Value = 32.63010831358289
$ws.Range("R5").Value = 293.670974822246
$ws.Range("S5").Value = 0.2203159030946636
$ws.Range("T5").Value = 0.2203159030946636
$ws.Range("I6").Value = 0.6840925621829359
$ws.Range("J6").Value = 0.684092562182936
$ws.Range("O6").Value = 0.5509544596378365
$ws.Range("P6").Value = 0.5509544596378364
$ws.Range("S6").Value = 0.3769038479397625
$ws.Range("T6").Value = 0.3769038479397625
$ws.Range("I7").Value = 0.6840925621829359
$ws.Range("J7").Value = 0.684092562182936
$ws.Range("O7").Value = 0.1269898489632735
$ws.Range("P7").Value = 0.1269898489632735
$ws.Range("S7").Value = 0.08687281114850984
$ws.Range("T7").Value = 0.08687281114850984
$ws.Range("G8").Value = 0.5566186666666667
$ws.Range("I8").Value = 0.1838612293955817
$ws.Range("J8").Value = 0.1838612293955817
$ws.Range("M8").Value = 15.75563966666667
$ws.Range("N8").Value = 47.266919
$ws.Range("O8").Value = 0.3220556913988901
$ws.Range("P8").Value = 0.32205569139889
$ws.Range("Q8").Value = 8.769883143740445
$ws.Range("R8").Value = 78.92894829366401
$ws.Range("S8").Value = 0.05921355535444401
$ws.Range("T8").Value = 0.05921355535444399
$ws.Range("G9").Value = 0.5566186666666667
$ws.Range("I9").Value = 0.1838612293955817
$ws.Range("J9").Value = 0.1838612293955817
$ws.Range("O9").Value = 0.5509544596378365
$ws.Range("P9").Value = 0.5509544596378364
$ws.Range("S9").Value = 0.101299164289991
$ws.Range("T9").Value = 0.101299164289991
$ws.Range("G10").Value = 0.5566186666666667
$ws.Range("I10").Value = 0.1838612293955817
$ws.Range("J10").Value = 0.1838612293955817
$ws.Range("O10").Value = 0.1269898489632735
$ws.Range("P10").Value = 0.1269898489632735
$ws.Range("Q10").Value = 3.458054509180445
$ws.Range("S10").Value = 0.02334850975114671
$ws.Range("T10").Value = 0.0233485097511467
